$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "26.443.15"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "1.611.35"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("E4").Value = "  -0.10%  "
Set-TextValue "D5" "211.74"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  +1.54%  "
Set-TextValue "D11" "0.0849"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").Value = "1.835.34"
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("D13").Value = "1.609.44"
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("E14").Value = "  -0.01%  "
Set-TextValue "D15" "0.508"
$ws.Range("E15").Value = "  -0.05%  "
Set-TextValue "D16" "63.56"
$ws.Range("E16").Value = "  -0.40%  "
Set-TextValue "D17" "234.66"
$ws.Range("E17").Value = "  +9.13%  "
$ws.Range("D18").Value = "26.436.12"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("E19").Value = "  +0.02%  "
Set-TextValue "D20" "7.64"
$ws.Range("E20").Value = "  +3.91%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  -0.19%  "
Set-TextValue "D23" "2.19"
$ws.Range("E23").Value = "  +4.11%  "
Set-TextValue "D24" "9.02"
$ws.Range("E24").Value = "  -0.22%  "
Set-TextValue "D25" "146.64"
$ws.Range("E25").Value = "  +1.37%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("E28").Value = "  +0.25%  "
Set-TextValue "D29" "15.45"
$ws.Range("E29").Value = "  +2.31%  "
$ws.Range("E30").Value = "  +0.94%  "
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("D32").Value = "1.492.93"
$ws.Range("E32").Value = "  +5.79%  "
$ws.Range("E33").Value = "  +1.30%  "
Set-TextValue "D34" "2.95"
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("E35").Value = "  -0.46%  "
$ws.Range("E36").Value = "  +1.34%  "
$ws.Range("E37").Value = "  -2.80%  "
$ws.Range("E38").Value = "  -0.12%  "
Set-TextValue "D39" "0.826"
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("E42").Value = "  +1.16%  "
Set-TextValue "D43" "0.928"
$ws.Range("E43").Value = "  -3.07%  "
$ws.Range("D44").Value = "1.748.23"
$ws.Range("E44").Value = "  +1.18%  "
Set-TextValue "D45" "0.761"
$ws.Range("E45").Value = "  -0.01%  "
Set-TextValue "D46" "61.26"
$ws.Range("E46").Value = "  +0.65%  "
Set-TextValue "D47" "89.84"
$ws.Range("E47").Value = "  +3.38%  "
$ws.Range("E48").Value = "  -2.18%  "
Set-TextValue "D49" "1.49"
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("E50").Value = "  -0.01%  "
Set-TextValue "D51" "0.0960"
$ws.Range("E51").Value = "  +0.95%  "
